$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Аркуш1" to "Sheet1"
$ws.Name = "Sheet1"

# --- Header row (A1:C1) ---
$ws.Range("A1").Value = "text"
$ws.Range("B1").Value = "agenda"
$ws.Range("C1").Value = "hi"

# --- Data rows, filled column-by-column (matches original authoring / shared-string order) ---
# Column A
$ws.Range("A2").Value = "спасибо"
$ws.Range("A3").Value = "сделать"
$ws.Range("A4").Value = "сказать"
$ws.Range("A5").Value = "сказать"
$ws.Range("A6").Value = "спасибо"
$ws.Range("A7").Value = "сказать"
$ws.Range("A8").Value = "сказать"
$ws.Range("A9").Value = "сказать"

# Column B
$ws.Range("B2").Value = "Благодарность"
$ws.Range("B3").Value = "Благодарность"
$ws.Range("B4").Value = "Благодарность"
$ws.Range("B5").Value = "Благодарность"
$ws.Range("B6").Value = "Благодарность"
$ws.Range("B7").Value = "Благодарность"
$ws.Range("B8").Value = "Благодарность"
$ws.Range("B9").Value = "Благодарность"

# Column C (numeric)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1

# --- Header formatting: bold font, thin box border, centered horizontal / top vertical alignment ---
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# --- Workbook default table/pivot style (matches source document metadata) ---
$wb.DefaultTableStyle = "TableStyleMedium9"
$wb.DefaultPivotTableStyle = "PivotStyleLight16"
